# IWEREYOU-44 Create function point files
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Update existing data rows (2-4) ----
$ws.Range("C2").Value2 = 33.15
$ws.Range("D2").Value2 = 1.5

$ws.Range("C3").Value2 = 51.85
$ws.Range("D3").Value2 = 8

$ws.Range("C4").Value2 = 28.05
$ws.Range("D4").Value2 = 1.5

# ---- New row 5: Manage Contacts (Create/ Show) ----
$ws.Range("B5").Value2 = "Manage Contacts (Create/ Show)"
$ws.Range("C5").Value2 = 35.7
$ws.Range("D5").Value2 = 5

# ---- New row 6: Registration + Login (computed) ----
$ws.Range("B6").Value2 = "Registration + Login"
$ws.Range("C6").Formula = "=SUM(C2:C3)"
$ws.Range("D6").Formula = "=SUM(D2:D3)"

# ---- New row 7: just a "Calculated Hours" label in D7 ----
$ws.Range("D7").Value2 = "Calculated Hours"

# ---- New row 8: Delete Contact ----
$ws.Range("B8").Value2 = "Delete Contact"
$ws.Range("C8").Value2 = 28.05
$ws.Range("D8").Formula = "=(C8-10.662)/7.2977"

# ---- New row 9: Invite a friend ----
$ws.Range("B9").Value2 = "Invite a friend"
$ws.Range("C9").Value2 = 47.6
$ws.Range("D9").Formula = "=(C9-10.662)/7.2977"

# ---- New row 10: View/ Update Profile ----
$ws.Range("B10").Value2 = "View/ Update Profile"
$ws.Range("C10").Value2 = 36.55
$ws.Range("D10").Formula = "=(C10-10.662)/7.2977"

# ---- New row 11: Manage Challenges ----
$ws.Range("B11").Value2 = "Manage Challenges"
$ws.Range("C11").Value2 = 44.2
$ws.Range("D11").Formula = "=(C11-10.662)/7.2977"

# ---- New row 12: Answer Challenge ----
$ws.Range("B12").Value2 = "Answer Challenge"
$ws.Range("C12").Value2 = 41.65
$ws.Range("D12").Formula = "=(C12-10.662)/7.2977"

# ---- Column widths ----
$ws.Columns.Item(2).ColumnWidth = 28.85546875
$ws.Columns.Item(4).ColumnWidth = 22.7109375

# ---- Selection ----
$ws.Range("C10").Select()
